$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 114, shifting existing rows 114-176 down to 115-177.
$ws.Rows("114:114").Insert()

# Populate the newly inserted row 114 with the new data point.
$ws.Cells.Item(114, 1).Value = 4
$ws.Cells.Item(114, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(114, 3).Value = "Los Lagos"
$ws.Cells.Item(114, 4).Value = 44518
$ws.Cells.Item(114, 4).NumberFormat = $ws.Cells.Item(115, 4).NumberFormat
$ws.Cells.Item(114, 5).Value = 10
$ws.Cells.Item(114, 6).Value = 100112044
$ws.Cells.Item(114, 7).Value = "Perejil"
$ws.Cells.Item(114, 8).Value = "Sin especificar"
$ws.Cells.Item(114, 9).Value = "Primera"
$ws.Cells.Item(114, 10).Value = 60
$ws.Cells.Item(114, 11).Value = 5000
$ws.Cells.Item(114, 12).Value = 5000
$ws.Cells.Item(114, 13).Value = 5000
$ws.Cells.Item(114, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(114, 15).Value = "Región Metropolitana"
$ws.Cells.Item(114, 16).Value = 1667
$ws.Cells.Item(114, 17).Value = 3
$ws.Cells.Item(114, 18).Value = "Hortaliza"
